$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A7 was a hardcoded literal (5); make it a formula matching the rest of the column
$ws.Range("A7").Formula = "=A6+1"

# Extend the weekly table down to rows 8 and 9 (two more weeks), following the
# exact same relative-formula pattern already used for rows 3-7
$ws.Range("A8").Formula = "=A7+1"
$ws.Range("B8").Formula = "=B7+7"
$ws.Range("C8").Formula = "=C7+7"

$ws.Range("A9").Formula = "=A8+1"
$ws.Range("B9").Formula = "=B8+7"
$ws.Range("C9").Formula = "=C8+7"

# Match the author's final selection (moved to A8 after the edit)
$ws.Range("A8").Select()
